# Applies the price / volume(1h) updates described by the commit diff.
# Column D ("Price") cells are stored as TEXT in the original workbook even
# when the value looks numeric (e.g. "579.38"), so values that would be
# auto-parsed as a number by Excel are written with a leading apostrophe to
# force a text (quote-prefixed) cell, matching the source formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Range, [string]$Value, [bool]$ForceText = $false)
    if ($ForceText) {
        $Range.Value = "'" + $Value
    } else {
        $Range.Value = $Value
    }
}

# Row 2
Set-TextCell $ws.Range("D2") '64.968.23' $false
Set-TextCell $ws.Range("E2") '  +0.17%  ' $false
# Row 3
Set-TextCell $ws.Range("D3") '3.147.51' $false
Set-TextCell $ws.Range("E3") '  -0.30%  ' $false
# Row 4
Set-TextCell $ws.Range("E4") '  +0.05%  ' $false
# Row 5
Set-TextCell $ws.Range("D5") '579.38' $true
Set-TextCell $ws.Range("E5") '  +1.48%  ' $false
# Row 6
Set-TextCell $ws.Range("D6") '148.48' $true
Set-TextCell $ws.Range("E6") '  -1.16%  ' $false
# Row 8
Set-TextCell $ws.Range("D8") '3.147.47' $false
Set-TextCell $ws.Range("E8") '  -0.29%  ' $false
# Row 9
Set-TextCell $ws.Range("E9") '  -0.72%  ' $false
# Row 10
Set-TextCell $ws.Range("E10") '  -2.74%  ' $false
# Row 11
Set-TextCell $ws.Range("E11") '  -1.10%  ' $false
# Row 12
Set-TextCell $ws.Range("E12") '  -1.40%  ' $false
# Row 13
Set-TextCell $ws.Range("D13") '0.0000263' $true
Set-TextCell $ws.Range("E13") '  +0.07%  ' $false
# Row 14
Set-TextCell $ws.Range("D14") '37.06' $true
Set-TextCell $ws.Range("E14") '  -3.67%  ' $false
# Row 15
Set-TextCell $ws.Range("D15") '3.664.06' $false
Set-TextCell $ws.Range("E15") '  -0.18%  ' $false
# Row 16
Set-TextCell $ws.Range("D16") '64.870.00' $false
Set-TextCell $ws.Range("E16") '  -0.09%  ' $false
# Row 17
Set-TextCell $ws.Range("D17") '3.161.84' $false
Set-TextCell $ws.Range("E17") '  +0.35%  ' $false
# Row 18
Set-TextCell $ws.Range("E18") '  -1.30%  ' $false
# Row 19
Set-TextCell $ws.Range("E19") '  +0.25%  ' $false
# Row 20
Set-TextCell $ws.Range("D20") '503.69' $true
Set-TextCell $ws.Range("E20") '  -2.39%  ' $false
# Row 21
Set-TextCell $ws.Range("D21") '15.06' $true
Set-TextCell $ws.Range("E21") '  +0.83%  ' $false
# Row 22
Set-TextCell $ws.Range("D22") '0.713' $true
Set-TextCell $ws.Range("E22") '  -3.47%  ' $false
# Row 23
Set-TextCell $ws.Range("D23") '15.13' $true
Set-TextCell $ws.Range("E23") '  -2.21%  ' $false
# Row 24
Set-TextCell $ws.Range("D24") '7.71' $true
Set-TextCell $ws.Range("E24") '  -1.88%  ' $false
# Row 25
Set-TextCell $ws.Range("D25") '84.19' $true
Set-TextCell $ws.Range("E25") '  -1.04%  ' $false
# Row 26
Set-TextCell $ws.Range("E26") '  +0.04%  ' $false
# Row 27
Set-TextCell $ws.Range("D27") '9.04' $true
Set-TextCell $ws.Range("E27") '  +0.99%  ' $false
# Row 28
Set-TextCell $ws.Range("D28") '2.91' $true
Set-TextCell $ws.Range("E28") '  -0.53%  ' $false
# Row 29
Set-TextCell $ws.Range("E29") '  -1.21%  ' $false
# Row 30
Set-TextCell $ws.Range("B30") 'Stacks' $false
Set-TextCell $ws.Range("C30") 'https://coinranking.com/coin/mMPrMcB7+stacks-stx' $false
Set-TextCell $ws.Range("D30") '2.77' $true
Set-TextCell $ws.Range("E30") '  +2.93%  ' $false
# Row 31
Set-TextCell $ws.Range("B31") 'EthereumClassic' $false
Set-TextCell $ws.Range("C31") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc' $false
Set-TextCell $ws.Range("D31") '27.48' $true
Set-TextCell $ws.Range("E31") '  -1.59%  ' $false
# Row 32
Set-TextCell $ws.Range("E32") '  +0.08%  ' $false
# Row 33
Set-TextCell $ws.Range("E33") '  -0.11%  ' $false
# Row 34
Set-TextCell $ws.Range("D34") '6.34' $true
Set-TextCell $ws.Range("E34") '  +2.32%  ' $false
# Row 35
Set-TextCell $ws.Range("D35") '6.47' $true
Set-TextCell $ws.Range("E35") '  -2.30%  ' $false
# Row 36
Set-TextCell $ws.Range("D36") '54.89' $true
Set-TextCell $ws.Range("E36") '  -1.54%  ' $false
# Row 37
Set-TextCell $ws.Range("D37") '0.0884' $true
Set-TextCell $ws.Range("E37") '  +1.79%  ' $false
# Row 38
Set-TextCell $ws.Range("D38") '474.04' $true
Set-TextCell $ws.Range("E38") '  -2.61%  ' $false
# Row 39
Set-TextCell $ws.Range("D39") '0.0413' $true
Set-TextCell $ws.Range("E39") '  -2.77%  ' $false
# Row 40
Set-TextCell $ws.Range("D40") '2.91' $true
Set-TextCell $ws.Range("E40") '  -3.25%  ' $false
# Row 41
Set-TextCell $ws.Range("D41") '8.72' $true
Set-TextCell $ws.Range("E41") '  +0.50%  ' $false
# Row 42
Set-TextCell $ws.Range("D42") '2.990.23' $false
Set-TextCell $ws.Range("E42") '  -4.24%  ' $false
# Row 43
Set-TextCell $ws.Range("D43") '0.116' $true
Set-TextCell $ws.Range("E43") '  -2.77%  ' $false
# Row 44
Set-TextCell $ws.Range("E44") '  -4.27%  ' $false
# Row 45
Set-TextCell $ws.Range("E45") '  -2.42%  ' $false
# Row 46
Set-TextCell $ws.Range("D46") '28.20' $true
Set-TextCell $ws.Range("E46") '  -5.21%  ' $false
# Row 47
Set-TextCell $ws.Range("D47") '0.0₃0589' $false
Set-TextCell $ws.Range("E47") '  +2.03%  ' $false
# Row 49
Set-TextCell $ws.Range("E49") '  -1.77%  ' $false
# Row 50
Set-TextCell $ws.Range("D50") '2.24' $true
Set-TextCell $ws.Range("E50") '  -2.80%  ' $false
# Row 51
Set-TextCell $ws.Range("E51") '  +14.08%  ' $false
